# Refactoring to use order instead of extensionOrder.
# This script updates the "Extensions_test" sheet (4th sheet) so that the
# EXTENSIONID column becomes a RELATION column that stores, for every other
# row, the EXTENSIONVALUE of the preceding row (i.e. a self relation instead
# of the old extensionId). It also fixes a duplicated "extensiontest13"
# value (row 14) to "extensiontest11", and updates the active sheet /
# selections on every worksheet to match the final, saved state of the
# workbook.

$wb = $excel.ActiveWorkbook

$wsCodeSchemes      = $wb.Worksheets.Item(1)
$wsCodes            = $wb.Worksheets.Item(2)
$wsExtensionSchemes = $wb.Worksheets.Item(3)
$wsExtensionsTest   = $wb.Worksheets.Item(4)

# --- Extensions_test sheet -------------------------------------------------

# Header: rename the EXTENSIONID column to RELATION.
$wsExtensionsTest.Range("D1").Value = "RELATION"

# Row 14 previously duplicated "extensiontest13"; it should reference the
# existing "extensiontest11" value instead.
$wsExtensionsTest.Range("A14").Value = "extensiontest11"

# Populate the new RELATION column: every other data row points back to the
# EXTENSIONVALUE of the row directly above it.
$wsExtensionsTest.Range("D3").Value  = "extensiontest1"
$wsExtensionsTest.Range("D5").Value  = "extensiontest3"
$wsExtensionsTest.Range("D7").Value  = "extensiontest5"
$wsExtensionsTest.Range("D9").Value  = "extensiontest7"
$wsExtensionsTest.Range("D11").Value = "extensiontest9"
$wsExtensionsTest.Range("D13").Value = "extensiontest11"
$wsExtensionsTest.Range("D15").Value = "extensiontest11"
$wsExtensionsTest.Range("D17").Value = "extensiontest15"
$wsExtensionsTest.Range("D19").Value = "extensiontest17"
$wsExtensionsTest.Range("D21").Value = "extensiontest19"
$wsExtensionsTest.Range("D23").Value = "extensiontest21"
$wsExtensionsTest.Range("D25").Value = "extensiontest23"

# --- Selections / active sheet ---------------------------------------------

$wsCodeSchemes.Activate()
$wsCodeSchemes.Range("A2").Select()

$wsExtensionSchemes.Activate()
$wsExtensionSchemes.Range("E2").Select()

# Extensions_test becomes the active (selected) sheet/tab.
$wsExtensionsTest.Activate()
$wsExtensionsTest.Range("D25").Select()
